$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "time" column (currently F) needs to move to the 3rd position (C),
# pushing systole/diastole/heartRate (currently C/D/E) one column to the
# right (to D/E/F). This is a cyclic left-rotation of the data held in
# columns C:F. We perform it as a sequence of direct-destination Cut
# operations (Range.Cut(Destination)) - using a helper column (H) to
# stage the original "time" column - because that preserves the exact
# original numeric values (bit-for-bit) and keeps the used range/
# dimension intact once the helper column is cleared again afterwards.

$rows = $ws.UsedRange.Rows.Count

$ws.Range("F1:F$rows").Cut($ws.Range("H1:H$rows")) | Out-Null
$ws.Range("E1:E$rows").Cut($ws.Range("F1:F$rows")) | Out-Null
$ws.Range("D1:D$rows").Cut($ws.Range("E1:E$rows")) | Out-Null
$ws.Range("C1:C$rows").Cut($ws.Range("D1:D$rows")) | Out-Null
$ws.Range("H1:H$rows").Cut($ws.Range("C1:C$rows")) | Out-Null

# Drop any residue left behind in the staging column so the sheet's
# dimension/used-range goes back to A:F.
$ws.Range("H1:H$rows").Clear() | Out-Null

# Reset the view: scroll back to the top-left and select B2 (matching
# the post-edit workbook state) instead of the old scrolled-down
# F92:F93 selection.
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("B2").Select() | Out-Null
